$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.765.76"
$ws.Range("E2").Value = "  +0.86%  "

# Row 3
$ws.Range("D3").Value = "1.648.42"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.50%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("E7").Value = "  +0.43%  "

# Row 8
$ws.Range("E8").Value = "  +0.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "

# Row 10
$ws.Range("E10").Value = "  +2.01%  "

# Row 12
$ws.Range("D12").Value = "1.879.05"
$ws.Range("E12").Value = "  +1.27%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.656.51"
$ws.Range("E13").Value = "  +1.69%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("E15").Value = "  +1.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("D17").Value = "26.767.04"
$ws.Range("E17").Value = "  +0.79%  "

# Row 18
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "

# Row 20
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +14.66%  "

# Row 22
$ws.Range("E22").Value = "  +1.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "

# Row 26
$ws.Range("E26").Value = "  +0.40%  "

# Row 27
$ws.Range("E27").Value = "  -0.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "

# Row 29
$ws.Range("E29").Value = "  +0.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "

# Row 31
$ws.Range("E31").Value = "  +1.66%  "

# Row 32
$ws.Range("E32").Value = "  -0.54%  "

# Row 33
$ws.Range("E33").Value = "  +0.90%  "

# Row 34
$ws.Range("D34").Value = "1.277.92"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.52%  "

# Row 36
$ws.Range("E36").Value = "  +2.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.67%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "

# Row 40
$ws.Range("E40").Value = "  +0.41%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.01%  "

# Row 42
$ws.Range("E42").Value = "  -0.98%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
$ws.Range("D44").Value = "1.790.68"
$ws.Range("E44").Value = "  +1.48%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.44%  "

# Row 47
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("E48").Value = "  +1.09%  "

# Row 49
$ws.Range("E49").Value = "  +0.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.42%  "

# Row 51
$ws.Range("E51").Value = "  +1.52%  "
